# LuanMenh.xlsx - "Tu vi dong cung kinh duong" edit
#
# Sheet2 holds a simple two-column (A/B) lookup list of "destiny" sentences.
# This change:
#   1) Rewords the three existing "Nguoi sinh nam ... co Tu Vi toa thu cung
#      Menh o Ngo" rows (18-20) to add the "khong gap Thien Hinh, Hoa Ky"
#      clause.
#   2) Appends 13 brand-new rows (21-33) with additional destiny sentences
#      (some duplicate existing phrases elsewhere in the sheet, some are new).
#   3) Leaves the scroll position at the top and moves the active selection
#      to the new last cell, B33.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Reword rows 18-20 ---------------------------------------------------
$ws.Cells.Item(18, 1).Value = "Sinh năm Giáp có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Hóa Kỵ"
$ws.Cells.Item(18, 2).Value = "Sinh năm Giáp có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Hóa Kỵ"

$ws.Cells.Item(19, 1).Value = "Sinh năm Đinh có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Hóa Kỵ"
$ws.Cells.Item(19, 2).Value = "Sinh năm Đinh có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Hóa Kỵ"

$ws.Cells.Item(20, 1).Value = "Sinh năm Kỷ có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Hóa Kỵ"
$ws.Cells.Item(20, 2).Value = "Sinh năm Kỷ có Tử Vi tọa thủ cung Mệnh ở Ngọ và không gặp Thiên Hình, Hóa Kỵ"

# --- New rows 21-33 --------------------------------------------------------
$ws.Cells.Item(21, 1).Value = "Nam mệnh sinh năm Nhâm có Tử Vi tọa thủ cung Mệnh ở Hợi"
$ws.Cells.Item(21, 2).Value = "Nam mệnh sinh năm Nhâm có Tử Vi tọa thủ cung Mệnh ở Hợi"

$ws.Cells.Item(22, 1).Value = "Nữ mệnh sinh năm Nhâm có Tử Vi tọa thủ cung Mệnh ở Dần"
$ws.Cells.Item(22, 2).Value = "Nữ mệnh sinh năm Nhâm có Tử Vi tọa thủ cung Mệnh ở Dần"

$ws.Cells.Item(23, 1).Value = "Nam mệnh sinh năm Giáp có Tử Vi tọa thủ cung Mệnh ở Hợi"
$ws.Cells.Item(23, 2).Value = "Nam mệnh sinh năm Giáp có Tử Vi tọa thủ cung Mệnh ở Hợi"

$ws.Cells.Item(24, 1).Value = "Nữ mệnh sinh năm Giáp có Tử Vi tọa thủ cung Mệnh ở Dần"
$ws.Cells.Item(24, 2).Value = "Nữ mệnh sinh năm Giáp có Tử Vi tọa thủ cung Mệnh ở Dần"

$ws.Cells.Item(25, 1).Value = "Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Dần"
$ws.Cells.Item(25, 2).Value = "Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Dần"

$ws.Cells.Item(26, 1).Value = "Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Thân"
$ws.Cells.Item(26, 2).Value = "Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Thân"

$ws.Cells.Item(27, 1).Value = "Tử Vi tọa thủ cung Mệnh đồng cung Thiên Phủ gặp Tả Phù, Hữu Bật"
$ws.Cells.Item(27, 2).Value = "Tử Vi tọa thủ cung Mệnh đồng cung Thiên Phủ gặp Tả Phù, Hữu Bật"

$ws.Cells.Item(28, 1).Value = "Tử Vi tọa thủ cung Mệnh gặp Tả Phù, Hữu Bật"
$ws.Cells.Item(28, 2).Value = "Tử Vi tọa thủ cung Mệnh gặp Tả Phù, Hữu Bật"

$ws.Cells.Item(29, 1).Value = "Thiên Phủ tọa thủ cung Mệnh gặp Tả Phù, Hữu Bật"
$ws.Cells.Item(29, 2).Value = "Thiên Phủ tọa thủ cung Mệnh gặp Tả Phù, Hữu Bật"

$ws.Cells.Item(30, 1).Value = "Tử Vi tọa thủ cung Mệnh đồng cung Kình Dương"
$ws.Cells.Item(30, 2).Value = "Tử Vi tọa thủ cung Mệnh đồng cung Kình Dương"

$ws.Cells.Item(31, 1).Value = "Thiên Phủ tọa thủ cung Mệnh đồng cung Kình Dương"
$ws.Cells.Item(31, 2).Value = "Thiên Phủ tọa thủ cung Mệnh đồng cung Kình Dương"

$ws.Cells.Item(32, 1).Value = "Tử Vi tọa thủ cung Mệnh và hội chiếu các sao: Thiên Phủ, Vũ Khúc, Thiên Tướng, Hóa Khoa, Hóa Lộc, Hóa Quyền, Long Trì, Phượng Các, Tả Phù, Hữu Bật, Quốc Ấn"
$ws.Cells.Item(32, 2).Value = "Tử Vi tọa thủ cung Mệnh và hội chiếu các sao: Thiên Phủ, Vũ Khúc, Thiên Tướng, Hóa Khoa, Hóa Lộc, Hóa Quyền, Long Trì, Phượng Các, Tả Phù, Hữu Bật, Quốc Ấn"

$ws.Cells.Item(33, 1).Value = "Tử Vi tọa thủ cung Mệnh và hội chiếu Địa Kiếp, Địa Không"
$ws.Cells.Item(33, 2).Value = "Tử Vi tọa thủ cung Mệnh và hội chiếu Địa Kiếp, Địa Không"

# --- Update the view: scroll back to top, select the new last cell --------
$ws.Activate()
$ws.Range("B33").Select()
